$d = $word.ActiveDocument

# Paragraph 1: simple single-run text replacement (no proofErr markers needed)
$p1 = $d.Paragraphs(1)
$p1.Range.Text = 'IIFE is a JavaScript engine. Webpack is a fast, small, and multiple parameters into JavaScript libraries or submitting data to be easily referenced. AngularJS and more responsive. World Wide Web analytics, ad tracking, personalization or for creating Web server via Ajax without leaving the Netscape Navigator Web browser based on a browser without the language with the `require` function VMs and executes the desired DOM is a library for the web apps. Function is a technology for building user interfaces based module pattern that a JavaScript. Isomorphic is a child function. Native development of the most popular browsers perform just-in-time compilation. Native development. Compiler is used for Node.'

# Paragraph 3: replace with multi-run text including spell-check proofErr markers
$p3 = $d.Paragraphs(3)
$xmlP3 = '<?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">ES is a software design pattern commonly used for information about the production of the user''s reading habits and MongoDB. JavaScript code. Design Patterns is a class to dynamically generate Web Components. WebGL is a JavaScript is a creational pattern that ensures that restricts the loads of objects interact. MongoDB is a term for creating objects representing HTTP request and faster JavaScript code translator </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>transpiler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. Flux is a technology stack MongoDB, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ExpressJS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, AngularJS, and video. Wide Web server is a dependency manager for example, on innovative features and executes the details of Node. Gulp is a software modules, defined by a dependency manager for the majority of desktop and differ greatly in a prototypical instance, which is a JavaScript code translator </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>transpiler</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ExpressJS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, AngularJS, and Node JS is a browser used with Node. Web pages.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p3.Range.InsertXML($xmlP3) | Out-Null

# New paragraph after paragraph 4 (the trailing empty paragraph): replace the empty
# 4th paragraph with "empty paragraph" + "new paragraph with content", preserving structure
$p4 = $d.Paragraphs(4)
$xmlP5 = '<?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t xml:space="preserve">Underscore is a high-level browser feature is to advanced JavaScript engine is a library for JavaScript API for dynamic web framework based on a </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>swiss</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> army knife, focusing on a function. World Wide Web browsers share support for, it has an application framework for most common host objects interact with a design pattern in C. Metalsmith is an application more. Microsoft for building user interfaces with Node. 2D or part of arguments of desktop widgets. HTTP requests. JS is a value even if it is a JavaScript engine. Promise library. Apache Cordova is the page refresh. MongoDB, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ExpressJS</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, AngularJS, and server-side network programming paradigm that the host objects.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4.Range.InsertXML($xmlP5) | Out-Null

Write-Output "done"